$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the weekly price-report values between row 2 and row 3
# (columns: D=Fecha, M=Volumen, N=Precio minimo, O=Precio maximo,
#  P=Precio promedio ponderado, S=Precio $/Kg)
$cols = @("D", "M", "N", "O", "P", "S")

foreach ($col in $cols) {
    $r2 = $ws.Range($col + "2")
    $r3 = $ws.Range($col + "3")
    $tmp = $r2.Value2
    $r2.Value2 = $r3.Value2
    $r3.Value2 = $tmp
}
